$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "ingrid.matos@mrv.com.br"
$ws.Range("B3").Value = "Bizagi"
$ws.Range("C3").Value = "Ferramenta de Planejamento"
$ws.Range("D3").Value = 5
$ws.Range("E3").Value = "Muito importante.`n"
